# Updates cryptos list - generated from OOXML diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.838.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.80%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.350.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.50%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.672"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.87%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'240.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.39%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'72.94"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.44%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.01%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.596"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.21%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.101"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.47%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'59.27"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +2.60%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'32.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +4.18%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'7.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.04%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +0.03%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.700.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.49%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'16.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -3.49%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.906"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.07%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.342.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.70%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'43.766.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.12%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -2.06%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.98%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'77.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.64%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'256.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.13%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +16.57%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.00%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'3.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.05%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.85%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'10.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.25%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "'EthereumClassic"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'22.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.30%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "'Monero"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'177.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.12%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "'Toncoin"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'2.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -5.38%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -0.65%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.70%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.0755"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.55%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'5.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -4.22%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'5.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.05%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.03%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -3.26%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'6.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -3.06%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +1.30%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'66.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +25.41%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +11.23%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'5.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +13.94%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.80%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'19.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.95%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.200"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.56%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -0.55%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'BinanceUSD"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.19%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'NEARProtocol"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'2.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.81%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'ARBITRUM"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'1.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.89%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'Aave"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'98.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.67%  "
$ws.Range("E51").Style = "Normal"
